{"js": "// Insert a new bold/italic run reading \"HELP ME\" right after the\n// \"1.2. Testing Schedule \" heading run, matching the target run's\n// size (14pt / w:sz=28) but keeping it as its own <w:r> element\n// (not merged into the neighboring run).\nconst body = context.document.body;\n\nconst results = body.search(\"1.2. Testing Schedule \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"1.2. Testing Schedule \" heading run.');\n}\n\nconst heading = results.items[0];\nconst inserted = heading.insertText(\"HELP ME\", \"After\");\ninserted.font.bold = true;\ninserted.font.italic = true;\n// Use a temporary distinct size so this new run is not silently\n// coalesced into the following (differently-anchored) run during\n// this same sync batch.\ninserted.font.size = 99;\nawait context.sync();\n\n// Fix the size up to match the heading (14pt) in its own sync batch,\n// after the run already exists as a separate element.\ninserted.font.size = 14;\nawait context.sync();\n", "ps1": "# Insert a new bold/italic run reading \"HELP ME\" right after the\n# \"1.2. Testing Schedule \" heading run, matching its size (14pt /\n# w:sz=28) but keeping it as its own <w:r> element (not merged into\n# the neighboring run).\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"1.2. Testing Schedule \"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$found = $rng.Find.Execute()\n\nif (-not $found) {\n    throw 'Could not find \"1.2. Testing Schedule \" heading run.'\n}\n\n$rng.Collapse(0)  # wdCollapseEnd - move insertion point to just after the found text\n$rng.InsertAfter(\"HELP ME\")\n$rng.Font.Bold = $true\n$rng.Font.Italic = $true\n# Use a temporary distinct size so this new run is not silently\n# coalesced into the following (differently-anchored) run.\n$rng.Font.Size = 99\n# Fix the size up to match the heading (14pt) as a separate write,\n# after the run already exists as a distinct element.\n$rng.Font.Size = 14\n"}
